$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.579.96"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.595.66"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.34%  "
$cell = $ws.Range("D5")
$cell.Value = "'207.85"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("E6").Value = "  -3.63%  "
$ws.Range("E7").Value = "  +0.37%  "
$cell = $ws.Range("D8")
$cell.Value = "'22.33"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -4.67%  "
$cell = $ws.Range("D10")
$cell.Value = "'0.0591"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("D12").Value = "1.823.04"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "1.580.96"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("E15").Value = "  -4.43%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "27.589.55"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell = $ws.Range("D17")
$cell.Value = "'63.32"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -3.11%  "
$cell = $ws.Range("D18")
$cell.Value = "'217.12"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -5.38%  "
$cell = $ws.Range("D19")
$cell.Value = "'7.38"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -4.05%  "
$ws.Range("D20").Value = "0.0₃0693"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("E21").Value = "  +0.35%  "
$cell = $ws.Range("D22")
$cell.Value = "'4.16"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -3.78%  "
$cell = $ws.Range("D23")
$cell.Value = "'9.59"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -4.67%  "
$ws.Range("E24").Value = "  -3.75%  "
$cell = $ws.Range("D25")
$cell.Value = "'152.48"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("E27").Value = "  +0.35%  "
$cell = $ws.Range("D28")
$cell.Value = "'15.08"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("E29").Value = "  -3.88%  "
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("E32").Value = "  -4.38%  "
$ws.Range("D33").Value = "1.376.42"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("E34").Value = "  -5.26%  "
$ws.Range("E35").Value = "  -3.99%  "
$ws.Range("E36").Value = "  -5.10%  "
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("E38").Value = "  -3.36%  "
$cell = $ws.Range("D39")
$cell.Value = "'0.540"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -2.91%  "
$cell = $ws.Range("D40")
$cell.Value = "'0.811"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +0.36%  "
$cell = $ws.Range("D42")
$cell.Value = "'0.973"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -4.21%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell = $ws.Range("D43")
$cell.Value = "'2.21"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Range("D44")
$cell.Value = "'5.36"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D45")
$cell.Value = "'1.78"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -2.88%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Range("D46")
$cell.Value = "'64.04"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").Value = "1.733.37"
$ws.Range("E47").Value = "  -2.21%  "
$cell = $ws.Range("D48")
$cell.Value = "'87.19"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("E50").Value = "  -4.25%  "
